# Update Data by bot, scripted by HH
# Applies the row-2 data refresh (report date 2020-06-30 -> 2018-12-31 plus
# the associated balance-sheet figures) described by the commit diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text / identifier columns -------------------------------------------------
# J2 (DATE_TYPE_CODE) looks numeric ("001"); prefix with an apostrophe so Excel
# keeps it as text instead of coercing it to the number 1, then restore the
# cell's original (unstyled) formatting so only the value itself changes.
$origStyle = $ws.Range("J2").Style
$ws.Range("J2").Value = "'001"
$ws.Range("J2").Style = $origStyle

# N2 (REPORT_DATE) is stored as plain text, not a real date value.
$ws.Range("N2").Value = "2018-12-31 00:00:00"

# --- Numeric columns -------------------------------------------------------
$ws.Range("O2").Value = 336551107.56
$ws.Range("P2").Value = 62631484.72
$ws.Range("Q2").Value = 10921947.43
$ws.Range("R2").Value = -47.7913803652
$ws.Range("S2").Value = 72701040.81
$ws.Range("T2").Value = 128.1333083934
$ws.Range("U2").Value = 27236101.22
$ws.Range("V2").Value = 63.8569780609
$ws.Range("W2").Value = 65732932.45
$ws.Range("X2").Value = 35879683.45
$ws.Range("Y2").Value = 79.0320065856

# Z2 / AA2 were empty before and now carry numeric values.
$ws.Range("Z2").Value = 663711.99
$ws.Range("AA2").Value = 68.2361768164

$ws.Range("AB2").Value = 270818175.11
$ws.Range("AC2").Value = 213.3747594815
$ws.Range("AD2").Value = 178.2772706778
$ws.Range("AE2").Value = 90.41422964039999
$ws.Range("AF2").Value = 358.6200252865
$ws.Range("AG2").Value = 19.5313374324
